# Apply the edit described by the diff:
# Insert a brand-new data row at row 537 (shifting the existing rows 537..625 down to 538..626),
# and populate the new row 537 with its data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 537; this shifts rows 537-625 down to 538-626
$ws.Rows(537).Insert()

# Populate the newly inserted row 537 with its values
$ws.Range("A537").Value = 3
$ws.Range("B537").Value = "Femacal de La Calera"
$ws.Range("C537").Value = "Coquimbo"
$ws.Range("D537").Value = 45180
$ws.Range("E537").Value = 5
$ws.Range("F537").Value = 100114013
$ws.Range("G537").Value = "Zanahoria"
$ws.Range("H537").Value = "Sin especificar"
$ws.Range("I537").Value = "Primera"
$ws.Range("J537").Value = 240
$ws.Range("K537").Value = 6500
$ws.Range("L537").Value = 7000
$ws.Range("M537").Value = 6750
$ws.Range("N537").Value = "`$/saco 20 kilos"
$ws.Range("O537").Value = "Provincia de Quillota"
$ws.Range("P537").Value = 338
$ws.Range("Q537").Value = 20
$ws.Range("R537").Value = "Hortaliza"
